$d = $word.ActiveDocument

$pairs = @(
    @("461×2=", "624×9="),
    @("643×6=", "523×7="),
    @("789×3=", "550×7="),
    @("115×9=", "466×7="),
    @("433×3=", "452×2="),
    @("308×8=", "238×2="),
    @("494×3=", "732×8="),
    @("470×3=", "486×5="),
    @("697×5=", "260×5="),
    @("177×4=", "278×9="),
    @("583×6=", "744×4="),
    @("598×8=", "822×6="),
    @("226×9=", "273×7="),
    @("359×5=", "539×8="),
    @("651×6=", "796×8="),
    @("657×6=", "843×2="),
    @("153×3=", "357×5="),
    @("957×2=", "852×6="),
    @("614×9=", "464×5="),
    @("714×7=", "727×4="),
    @("724×5=", "698×6="),
    @("951×6=", "526×7="),
    @("589×3=", "725×5="),
    @("116×9=", "993×6="),
    @("260×6=", "748×2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
